# Apply cryptos list update (Tue Jul 30 15:42:27 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'66.126.23"
$ws.Range('E2').Value = '  -2.65%  '
$ws.Range('D3').Value = "'3.317.59"
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = "'573.41"
$ws.Range('E5').Value = '  -1.53%  '
$ws.Range('D6').Value = "'181.07"
$ws.Range('E6').Value = '  -2.84%  '
$ws.Range('D7').Value = "'0.617"
$ws.Range('E7').Value = '  +3.57%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = "'0.127"
$ws.Range('E9').Value = '  -1.62%  '
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('E11').Value = '  -1.52%  '
$ws.Range('D12').Value = "'3.896.11"
$ws.Range('E12').Value = '  +0.60%  '
$ws.Range('E13').Value = '  -1.06%  '
$ws.Range('D14').Value = "'26.67"
$ws.Range('E14').Value = '  -2.92%  '
$ws.Range('D15').Value = "'66.261.25"
$ws.Range('E15').Value = '  -2.65%  '
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('D17').Value = "'3.343.30"
$ws.Range('E17').Value = '  +0.95%  '
$ws.Range('D18').Value = "'433.45"
$ws.Range('E18').Value = '  -3.12%  '
$ws.Range('D19').Value = "'13.54"
$ws.Range('D20').Value = "'5.64"
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('D21').Value = "'7.54"
$ws.Range('E21').Value = '  -2.33%  '
$ws.Range('D22').Value = "'73.22"
$ws.Range('E22').Value = '  -2.28%  '
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').Value = "'0.520"
$ws.Range('E24').Value = '  +1.14%  '
$ws.Range('D25').Value = "'3.471.90"
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('E26').Value = '  -2.28%  '
$ws.Range('E27').Value = '  +3.02%  '
$ws.Range('D28').Value = "'9.04"
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('D29').Value = "'0.999"
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('D30').Value = "'1.94"
$ws.Range('E30').Value = '  -1.69%  '
$ws.Range('D31').Value = "'22.68"
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').Value = "'1.00"
$ws.Range('D33').Value = "'5.22"
$ws.Range('E33').Value = '  -2.36%  '
$ws.Range('D34').Value = "'6.75"
$ws.Range('E34').Value = '  -0.60%  '
$ws.Range('D35').Value = "'1.21"
$ws.Range('E35').Value = '  -2.75%  '
$ws.Range('E36').Value = '  -2.86%  '
$ws.Range('D37').Value = "'159.60"
$ws.Range('E37').Value = '  -2.43%  '
$ws.Range('D38').Value = "'27.69"
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = "'2.845.65"
$ws.Range('E39').Value = '  +5.28%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = "'1.79"
$ws.Range('E40').Value = '  -3.72%  '
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('D42').Value = "'4.42"
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('D43').Value = "'6.17"
$ws.Range('E43').Value = '  -3.07%  '
$ws.Range('D44').Value = "'40.52"
$ws.Range('E44').Value = '  -0.46%  '
$ws.Range('D45').Value = "'0.0665"
$ws.Range('E45').Value = '  -1.06%  '
$ws.Range('E46').Value = '  -2.10%  '
$ws.Range('D47').Value = "'24.11"
$ws.Range('E47').Value = '  -2.21%  '
$ws.Range('D48').Value = "'325.64"
$ws.Range('E48').Value = '  -0.12%  '
$ws.Range('D49').Value = "'0.0271"
$ws.Range('E49').Value = '  -1.60%  '
$ws.Range('E50').Value = '  +1.78%  '
$ws.Range('D51').Value = "'0.972"
$ws.Range('E51').Value = '  -1.66%  '
